$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price values are stored as literal text (matches source formatting)
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D17", "D19", "D20", "D22", "D23", "D25", "D26", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D38", "D39", "D40", "D41", "D43", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "47.306.71"
$ws.Range("E2").Value = "  +4.90%  "
$ws.Range("D3").Value = "2.507.47"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "323.99"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "105.48"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").Value = "37.05"
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "18.48"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "7.25"
$ws.Range("E14").Value = "  +4.18%  "
$ws.Range("D15").Value = "2.907.71"
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").Value = "2.544.69"
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "47.298.53"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").Value = "0.0₃0940"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Value = "70.97"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("D23").Value = "252.77"
$ws.Range("E23").Value = "  +3.72%  "
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "26.35"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +5.41%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "35.39"
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("D31").Value = "0.135"
$ws.Range("E31").Value = "  +5.70%  "
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").Value = "19.64"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "0.0779"
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("D39").Value = "2.97"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "121.79"
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("D43").Value = "21.75"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").Value = "1.980.49"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.80"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.15"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "5.43"
$ws.Range("E50").Value = "  +16.22%  "
$ws.Range("D51").Value = "79.64"
$ws.Range("E51").Value = "  +4.58%  "
